$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12, column A was stored as text in the source file; the new export
# writes it as a genuine number, so normalize it here.
$ws.Range("A12").Value = 76442780

# Append the new payment row (row 13): phone, amount, method, timestamp,
# original_amount, discount_applied, final_amount.

# "A13" looks numeric but, like the other phone numbers already captured
# as strings in this sheet, must stay text. Force text via NumberFormat,
# then drop the format back to Normal so no stray number format lingers
# on the cell itself.
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "79172233"
$ws.Range("A13").Style = "Normal"

# "B13" has no value for this row (same as B12), but the cell itself is
# still present in the sheet, so write it and then touch formatting
# (back to its default) so the now-blank cell isn't dropped entirely.
$ws.Range("B13").Value = ""
$ws.Range("B13").Font.Bold = $false

$ws.Range("C13").Value = "Credit Card"
$ws.Range("D13").Value = "2025-08-18T08:33:09"
$ws.Range("E13").Value = 30
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 30
